$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 18: new logboek entry (Profile Pagina & geboden.html / rating as float)
# ---------------------------------------------------------------------------

$text = "Profile Pagina & geboden.html: Als klusjesman nog geen ratings heeft, error bij ophalen van pagina -> gefixt via try catch met (NA/5). Ook overal rating op 10 veranderd naar rating op 5 (geboden,info,review,...). [maincontroller]: code in maincontroller was aangepast voor ophalen van toegewezen klusjes via de OfferRepo, maar daardoor waren er bugs: Elke klusjesman die ooit geboden heeft op een task zag deze verschijnen op hun profiel in alle volgende fasen ookal hadden ze deze niet toegewezen gekregen. Zij konden deze ook markeren als uitgevoerd etc. Daarnaast waren er model variabelen aangepast op 1 plek maar niet op de andere waardoor de volledig afgewerkte klusjes niet meer verschenen -> gefixt. geboden,review,.. website logo in de menubalk linksboven anders gelinkt in general.html zodat deze zichtbaar is op deze pagina's. ALLES code in commentaar weghalen in elke klasse/controller/... Index & Profile: GEBODEN task weergeeft aantal biedingen rechtsonder op postit. Enkel eigenaar kan info van klusjesmannen zien op toekenningspagina.  ALLES nog is getest van start to finish + opgave nagelezen om te kijken of alles erin zit. Task deleten als er al geboden is nog toegevoegd. RATING: Aangepast zodat het als float wordt berekend en weergegeven."

$cell = $ws.Range("A18")
$cell.Value = $text

# Apply explicit font + bold/italic runs to match the original logboek styling
# (every run keeps Calibri/12/automatic-black like the rest of the log, with
# bold used for the labels and italic for the "[maincontroller]:" note)
$cell.Characters(1,29).Font.Name = "Calibri"
$cell.Characters(1,29).Font.Size = 12
$cell.Characters(1,29).Font.ColorIndex = 1
$cell.Characters(1,29).Font.Bold = $true

$cell.Characters(30,158).Font.Name = "Calibri"
$cell.Characters(30,158).Font.Size = 12
$cell.Characters(30,158).Font.ColorIndex = 1

$cell.Characters(188,23).Font.Name = "Calibri"
$cell.Characters(188,23).Font.Size = 12
$cell.Characters(188,23).Font.ColorIndex = 1
$cell.Characters(188,23).Font.Bold = $true

$cell.Characters(211,3).Font.Name = "Calibri"
$cell.Characters(211,3).Font.Size = 12
$cell.Characters(211,3).Font.ColorIndex = 1

$cell.Characters(214,17).Font.Name = "Calibri"
$cell.Characters(214,17).Font.Size = 12
$cell.Characters(214,17).Font.ColorIndex = 1
$cell.Characters(214,17).Font.Italic = $true

$cell.Characters(231,478).Font.Name = "Calibri"
$cell.Characters(231,478).Font.Size = 12
$cell.Characters(231,478).Font.ColorIndex = 1

$cell.Characters(709,17).Font.Name = "Calibri"
$cell.Characters(709,17).Font.Size = 12
$cell.Characters(709,17).Font.ColorIndex = 1
$cell.Characters(709,17).Font.Bold = $true

$cell.Characters(726,58).Font.Name = "Calibri"
$cell.Characters(726,58).Font.Size = 12
$cell.Characters(726,58).Font.ColorIndex = 1

$cell.Characters(784,12).Font.Name = "Calibri"
$cell.Characters(784,12).Font.Size = 12
$cell.Characters(784,12).Font.ColorIndex = 1
$cell.Characters(784,12).Font.Bold = $true

$cell.Characters(796,43).Font.Name = "Calibri"
$cell.Characters(796,43).Font.Size = 12
$cell.Characters(796,43).Font.ColorIndex = 1

$cell.Characters(839,5).Font.Name = "Calibri"
$cell.Characters(839,5).Font.Size = 12
$cell.Characters(839,5).Font.ColorIndex = 1
$cell.Characters(839,5).Font.Bold = $true

$cell.Characters(844,59).Font.Name = "Calibri"
$cell.Characters(844,59).Font.Size = 12
$cell.Characters(844,59).Font.ColorIndex = 1

$cell.Characters(903,15).Font.Name = "Calibri"
$cell.Characters(903,15).Font.Size = 12
$cell.Characters(903,15).Font.ColorIndex = 1
$cell.Characters(903,15).Font.Bold = $true

$cell.Characters(918,135).Font.Name = "Calibri"
$cell.Characters(918,135).Font.Size = 12
$cell.Characters(918,135).Font.ColorIndex = 1

$cell.Characters(1053,5).Font.Name = "Calibri"
$cell.Characters(1053,5).Font.Size = 12
$cell.Characters(1053,5).Font.ColorIndex = 1
$cell.Characters(1053,5).Font.Bold = $true

$cell.Characters(1058,136).Font.Name = "Calibri"
$cell.Characters(1058,136).Font.Size = 12
$cell.Characters(1058,136).Font.ColorIndex = 1

$cell.Characters(1194,6).Font.Name = "Calibri"
$cell.Characters(1194,6).Font.Size = 12
$cell.Characters(1194,6).Font.ColorIndex = 1
$cell.Characters(1194,6).Font.Bold = $true

$cell.Characters(1200,62).Font.Name = "Calibri"
$cell.Characters(1200,62).Font.Size = 12
$cell.Characters(1200,62).Font.ColorIndex = 1

# Date of the entry
$ws.Range("B18").Value = (Get-Date -Year 2024 -Month 12 -Day 23 -Hour 0 -Minute 0 -Second 0)

# Hours worked - shown with the new "right aligned, one decimal" number style
$hoursCell = $ws.Range("C18")
$hoursCell.Value = 6
$hoursCell.NumberFormat = "0.0"
$hoursCell.HorizontalAlignment = -4152

# ---------------------------------------------------------------------------
# Update the active selection to reflect where the author ended up (row 18)
# ---------------------------------------------------------------------------
$ws.Range("E18").Select()

Write-Output "done"
